# Add I0 and IF columns (I and J) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers in row 1, matching the style of the existing header cells (e.g. H1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-18 for columns I (I0) and J (IF).
$values = @{
    2  = @(6, 6)
    3  = @(8, 8)
    4  = @(5, 6)
    5  = @(4, 4)
    6  = @(7, 7)
    7  = @(7, 8)
    8  = @(7, 7)
    9  = @(7, 7)
    10 = @(9, 9)
    11 = @(8, 8)
    12 = @(8, 8)
    13 = @(7, 7)
    14 = @(6, 6)
    15 = @(5, 6)
    16 = @(6, 6)
    17 = @(8, 8)
    18 = @(6, 6)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
